# Add language selection on headers
# - Rename sheet "general" -> "General"
# - Column C holds the Japanese (lang_jp) translations; give every data
#   row in column C its own localized text (previously it just mirrored
#   the English text).
# - Columns B, D, E and F all now show the same (English) text instead of
#   each column holding a distinct value.
# - Japanese cells (column C, rows with data) get a distinct look: a new
#   font ("Inherit", color #212121), a white solid fill, left aligned and
#   wrapped text.

$wb = $excel.ActiveWorkbook

# ---- rename first sheet ----
$wsGeneral = $wb.Worksheets.Item(1)
$wsGeneral.Name = "General"

$wsLogin = $wb.Worksheets.Item(2)

function Set-JapaneseCellStyle($cell) {
    $cell.Interior.Color = 16777215
    $cell.Font.Name = "Inherit"
    $cell.Font.Color = 2171169
    $cell.HorizontalAlignment = -4131
    $cell.WrapText = $true
}

# ================= General sheet =================
# Row 2: txtCopyright
$wsGeneral.Range("B2").Value = "Copyright Blockpass ©2018"
$wsGeneral.Range("C2").Value = "著作権 Blockpass ©2018"
$wsGeneral.Range("D2").Value = "Copyright Blockpass ©2018"
$wsGeneral.Range("E2").Value = "Copyright Blockpass ©2018"
$wsGeneral.Range("F2").Value = "Copyright Blockpass ©2018"
Set-JapaneseCellStyle $wsGeneral.Range("C2")

# ================= Login sheet =================
# Row 2: txtLogin
$wsLogin.Range("B2").Value = "Login"
$wsLogin.Range("C2").Value = "ログイン"
$wsLogin.Range("D2").Value = "Login"
$wsLogin.Range("E2").Value = "Login"
$wsLogin.Range("F2").Value = "Login"
Set-JapaneseCellStyle $wsLogin.Range("C2")

# Row 3: txtUsernamePlaceholder
$wsLogin.Range("B3").Value = "username"
$wsLogin.Range("C3").Value = "ユーザー名"
$wsLogin.Range("D3").Value = "username"
$wsLogin.Range("E3").Value = "username"
$wsLogin.Range("F3").Value = "username"
Set-JapaneseCellStyle $wsLogin.Range("C3")

# Row 4: txtInvalidUsername
$wsLogin.Range("B4").Value = "Invalid user name"
$wsLogin.Range("C4").Value = "無効なユーザー名"
$wsLogin.Range("D4").Value = "Invalid user name"
$wsLogin.Range("E4").Value = "Invalid user name"
$wsLogin.Range("F4").Value = "Invalid user name"
Set-JapaneseCellStyle $wsLogin.Range("C4")

# Row 5: txtPasswordPlaceholder
$wsLogin.Range("B5").Value = "password"
$wsLogin.Range("C5").Value = "パスワード"
$wsLogin.Range("D5").Value = "password"
$wsLogin.Range("E5").Value = "password"
$wsLogin.Range("F5").Value = "password"
Set-JapaneseCellStyle $wsLogin.Range("C5")

# Row 6: txtInvalidPassword
$wsLogin.Range("B6").Value = "Invalid password"
$wsLogin.Range("C6").Value = "無効なパスワード"
$wsLogin.Range("D6").Value = "Invalid password"
$wsLogin.Range("E6").Value = "Invalid password"
$wsLogin.Range("F6").Value = "Invalid password"
Set-JapaneseCellStyle $wsLogin.Range("C6")

Write-Host "Applied language selection changes"
